# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report record is inserted as row 17 of the data table,
# pushing all subsequent records (old rows 17-53) down by one row
# (new rows 18-54). The sheet's used range grows from A1:R53 to A1:R54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 - shifts rows 17..53 down to 18..54
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record
$ws.Cells.Item(17, 1).Value  = 7
$ws.Cells.Item(17, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value  = "Ñuble"
$ws.Cells.Item(17, 4).Value  = 44838
$ws.Cells.Item(17, 5).Value  = 16
$ws.Cells.Item(17, 6).Value  = 100112001
$ws.Cells.Item(17, 7).Value  = "Berenjena"
$ws.Cells.Item(17, 8).Value  = "Sin especificar"
$ws.Cells.Item(17, 9).Value  = "Primera"
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 12000
$ws.Cells.Item(17, 12).Value = 13000
$ws.Cells.Item(17, 13).Value = 12500
$ws.Cells.Item(17, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(17, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value = 208
$ws.Cells.Item(17, 17).Value = 60
$ws.Cells.Item(17, 18).Value = "Hortaliza"
